# "fixed excelfile" — the Process-Commodity sheet had a single row whose
# Process/Commodity names had been mashed together into one shared-string
# ("WaterplantElec") with the Commodity column left blank. Split it back
# into the correct two columns: Process="Waterplant", Commodity="Elec"
# (row 7: Waterplant / Elec / out / 1), matching every other row's layout
# in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Process-Commodity")

$ws.Range("A7").Value = "Waterplant"
$ws.Range("B7").Value = "Elec"

# Leave the sheet focused on the row we just fixed, and make
# "Process-Commodity" the active tab/sheet (it was "Process" before).
$ws.Activate() | Out-Null
$ws.Range("B8").Select() | Out-Null
